$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("W2").Value = "[1002003004001,1002003004002]"
$ws.Range("W3").Value = "[1002003004001,1002003004002]"

$ws.Range("W2:W3").Font.Name = "Calibri"
$ws.Range("W2:W3").Font.Size = 11
$ws.Range("W2:W3").WrapText = $false
